$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.059.57"
$ws.Range("E2").Value = "  -3.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.964.02"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.87"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.34"
$ws.Range("E6").Value = "  +4.49%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.511"
$ws.Range("E8").Value = "  +3.19%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.953.03"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  -2.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.83"
$ws.Range("E11").Value = "  -5.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  +3.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000221"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.80"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.451.47"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.953.48"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.68"
$ws.Range("E18").Value = "  +9.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.792.89"
$ws.Range("E19").Value = "  -4.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "416.65"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.12"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("E22").Value = "  +4.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.97"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.92"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.37"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.53"
$ws.Range("E29").Value = "  +6.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.98"
$ws.Range("E30").Value = "  +6.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.104"
$ws.Range("E32").Value = "  +12.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.10"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.61"
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.934"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.09"
$ws.Range("E36").Value = "  -6.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.32"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0660"
$ws.Range("E38").Value = "  +2.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.32"
$ws.Range("E39").Value = "  +6.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.56"
$ws.Range("E40").Value = "  +8.74%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0347"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "376.30"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.648.45"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.237"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.74"
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.38"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("E51").Value = "  +2.09%  "
